$wb = $excel.ActiveWorkbook

# "Flight Mission Cycle" sheet: append a new "Typing" setting row (Duration 60)
$ws2 = $wb.Worksheets.Item("Flight Mission Cycle")
$ws2.Range("A5").Value = "Typing"
$ws2.Range("B5").Value = 60

# "Settings" sheet: move the stored selection to A5
$ws3 = $wb.Worksheets.Item("Settings")
[void]$ws3.Range("A5").Select()

# Re-activate "Flight Mission Cycle" (keeps it the tab-selected sheet) and
# move its stored selection to M19
[void]$ws2.Activate()
[void]$ws2.Range("M19").Select()
